$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fix the ID value for the "Engine cooling unit" system row (was a typo: EN_A0006 -> EN_A0600)
$ws.Range("G2").Value = "EN_A0600"

# Update the active selection to G3
$ws.Activate()
$ws.Range("G3").Select()
